$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove trailing whitespace from the RefSeq Genome Accession of
# "Leptospira iowaensis Hamond et al. 2025" (row 85, column E)
$ws.Range("E85").Value = "GCF_040833965.1"

# Classify the five newly-added n=5 species with their P1/P2/S1 status
# (column I) for the latest batch of Leptospira species (rows 83-87)
$ws.Range("I83").Value = "P2"
$ws.Range("I84").Value = "P1"
$ws.Range("I85").Value = "S1"
$ws.Range("I86").Value = "S1"
$ws.Range("I87").Value = "S1"

# Update the view's scroll position / active selection to match where the
# editor ended up working
$excel.ActiveWindow.ScrollRow = 75
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("I89").Select()

$wb.Save()
